# Update "想去人数" (number of people interested) counts in the
# "展览" (exhibitions) and "全部类型" (all types) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 27
$ws1.Range("F4").Value = 53
$ws1.Range("F5").Value = 41
$ws1.Range("F6").Value = 261
$ws1.Range("F7").Value = 3686
$ws1.Range("F9").Value = 4304
$ws1.Range("F11").Value = 1073
$ws1.Range("F12").Value = 54

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 27
$ws4.Range("F4").Value = 53
$ws4.Range("F5").Value = 41
$ws4.Range("F7").Value = 261
$ws4.Range("F8").Value = 3686
$ws4.Range("F10").Value = 4304
$ws4.Range("F12").Value = 1073
$ws4.Range("F13").Value = 54
